$d = $word.ActiveDocument

$replacements = @(
    @{ old = "822×9=7398"; new = "999×2=1998" },
    @{ old = "417×8=3336"; new = "803×6=4818" },
    @{ old = "761×3=2283"; new = "368×2=736" },
    @{ old = "357×8=2856"; new = "741×3=2223" },
    @{ old = "666×9=5994"; new = "925×3=2775" },
    @{ old = "554×9=4986"; new = "116×3=348" },
    @{ old = "902×5=4510"; new = "644×6=3864" },
    @{ old = "599×7=4193"; new = "180×3=540" },
    @{ old = "977×9=8793"; new = "453×8=3624" },
    @{ old = "910×6=5460"; new = "528×7=3696" },
    @{ old = "836×7=5852"; new = "814×7=5698" },
    @{ old = "831×3=2493"; new = "158×3=474" },
    @{ old = "537×8=4296"; new = "481×8=3848" },
    @{ old = "889×2=1778"; new = "598×2=1196" },
    @{ old = "925×4=3700"; new = "378×7=2646" },
    @{ old = "965×2=1930"; new = "128×4=512" },
    @{ old = "870×7=6090"; new = "903×5=4515" },
    @{ old = "310×2=620"; new = "804×8=6432" },
    @{ old = "310×6=1860"; new = "849×3=2547" },
    @{ old = "676×9=6084"; new = "473×3=1419" },
    @{ old = "521×4=2084"; new = "110×3=330" },
    @{ old = "781×6=4686"; new = "686×8=5488" },
    @{ old = "596×3=1788"; new = "590×3=1770" },
    @{ old = "236×2=472"; new = "350×7=2450" },
    @{ old = "403×8=3224"; new = "684×4=2736" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
